$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.142.16'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.656.34'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  -0.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.70'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5245'
$ws.Range('E6').Value = '  -1.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2621'
$ws.Range('E8').Value = '  -0.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06294'
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.58'
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07809'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.497'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.661.74'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.884.32'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5551'
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8002'
$ws.Range('E16').Value = '  -2.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.06'
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.154.73'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.639'
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '195.56'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.964'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('E24').Value = '  -0.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.68'
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('E26').Value = '  -1.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.166'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  -0.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.498'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05719'
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.272'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.490'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.353'
$ws.Range('E33').Value = '  +2.31%  '
$ws.Range('E34').Value = '  -1.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.803'
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9528'
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.417'
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01597'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('E40').Value = '  +2.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.058.66'
$ws.Range('E41').Value = '  +1.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8444'
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.004'
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.48'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.794.74'
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('E46').Value = '  +0.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05395'
$ws.Range('E47').Value = '  +4.51%  '
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4399'
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.018'
$ws.Range('E51').Value = '  -0.54%  '
